# ADD: new intent and text to answer
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (인사 / 안녕하세요): drop the answer-image hyperlink ---
$ws.Range("E2").Hyperlinks.Delete()
$ws.Range("E2").ClearContents()
$ws.Rows.Item(2).RowHeight = 51.75

# --- Row 3: "인사"/반가워요 becomes new "만남" intent with a warning answer ---
$ws.Range("A3").Value = "만남"
$ws.Range("D3").Value = "[경고]`n상대방이 갑작스러운 만남을`n요구할 경우에는 항상 조심하세요!"
$ws.Range("E3").Hyperlinks.Delete()
$ws.Range("E3").ClearContents()
$ws.Rows.Item(3).RowHeight = 69

# --- Row 4: old 주문/B_FOOD row becomes new "금전" intent with a warning answer ---
$ws.Range("A4").Value = "금전"
$ws.Range("B4").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("D4").Value = "[경고]`n상대방이 갑작스러운 금전적`n요구할 경우에는 항상 조심하세요!"
$ws.Rows.Item(4).RowHeight = 69

# --- Row 5: keep "주문" intent label, clear the rest ---
$ws.Range("A5").Value = "주문"
$ws.Range("B5").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Rows.Item(5).RowHeight = 17.25

# --- Row 6: keep "예약" intent label, clear the rest ---
$ws.Range("A6").Value = "예약"
$ws.Range("B6").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Rows.Item(6).RowHeight = 17.25

# --- Row 7: 욕설 gets a new warning-style answer text ---
$ws.Range("A7").Value = "욕설"
$ws.Range("D7").Value = "[경고]`n상대방과 나를 위해 욕설 사용은 자제해주세요."
$ws.Rows.Item(7).RowHeight = 51.75

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("D7").Select()
